$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (B,C,D,E) keep their original Text storage type,
# matching the source inlineStr cells (avoids Excel auto-coercing
# numeric-looking strings like "1.00" or "9.20" into numbers).
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '63.646.84'
$ws.Range("E2").Value = '  +0.28%  '

# Row 3
$ws.Range("D3").Value = '2.619.47'
$ws.Range("E3").Value = '  -0.58%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '594.65'
$ws.Range("E5").Value = '  -1.59%  '

# Row 6
$ws.Range("D6").Value = '150.07'
$ws.Range("E6").Value = '  +2.09%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  +0.27%  '

# Row 10
$ws.Range("D10").Value = '5.69'
$ws.Range("E10").Value = '  +1.81%  '

# Row 11
$ws.Range("E11").Value = '  +3.09%  '

# Row 12
$ws.Range("E12").Value = '  -1.05%  '

# Row 13
$ws.Range("D13").Value = '27.68'
$ws.Range("E13").Value = '  +0.44%  '

# Row 14
$ws.Range("D14").Value = '3.092.99'
$ws.Range("E14").Value = '  -0.58%  '

# Row 15
$ws.Range("D15").Value = '63.497.88'
$ws.Range("E15").Value = '  +0.27%  '

# Row 16
$ws.Range("D16").Value = '0.0000151'
$ws.Range("E16").Value = '  +2.45%  '

# Row 17
$ws.Range("D17").Value = '2.586.63'
$ws.Range("E17").Value = '  -2.34%  '

# Row 18
$ws.Range("D18").Value = '12.34'
$ws.Range("E18").Value = '  +6.88%  '

# Row 19
$ws.Range("E19").Value = '  +1.85%  '

# Row 20
$ws.Range("D20").Value = '347.21'
$ws.Range("E20").Value = '  +0.87%  '

# Row 21
$ws.Range("D21").Value = '6.86'
$ws.Range("E21").Value = '  -0.68%  '

# Row 22
$ws.Range("E22").Value = '  -0.26%  '

# Row 23
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  +3.07%  '

# Row 24
$ws.Range("D24").Value = '66.43'
$ws.Range("E24").Value = '  -0.59%  '

# Row 25
$ws.Range("E25").Value = '  +9.81%  '

# Row 26
$ws.Range("D26").Value = '9.20'
$ws.Range("E26").Value = '  +1.83%  '

# Row 27
$ws.Range("D27").Value = '1.66'
$ws.Range("E27").Value = '  -2.06%  '

# Row 28
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '8.22'
$ws.Range("E28").Value = '  +2.81%  '

# Row 29
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '547.94'
$ws.Range("E29").Value = '  -2.40%  '

# Row 30
$ws.Range("E30").Value = '  -0.88%  '

# Row 31
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.10%  '

# Row 32
$ws.Range("E32").Value = '  -1.33%  '

# Row 33
$ws.Range("D33").Value = '0.0₃0843'
$ws.Range("E33").Value = '  -1.20%  '

# Row 34
$ws.Range("D34").Value = '1.74'
$ws.Range("E34").Value = '  -0.35%  '

# Row 35
$ws.Range("D35").Value = '5.24'
$ws.Range("E35").Value = '  +0.63%  '

# Row 36
$ws.Range("D36").Value = '168.38'
$ws.Range("E36").Value = '  +0.65%  '

# Row 37
$ws.Range("D37").Value = '0.413'
$ws.Range("E37").Value = '  +1.71%  '

# Row 38
$ws.Range("E38").Value = '  +0.03%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.94'
$ws.Range("E39").Value = '  -0.43%  '

# Row 40
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D40").Value = '19.40'
$ws.Range("E40").Value = '  +1.27%  '

# Row 41
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("D42").Value = '166.51'
$ws.Range("E42").Value = '  -0.58%  '

# Row 43
$ws.Range("D43").Value = '39.83'
$ws.Range("E43").Value = '  -0.56%  '

# Row 44
$ws.Range("D44").Value = '3.91'
$ws.Range("E44").Value = '  +3.32%  '

# Row 45
$ws.Range("D45").Value = '0.0591'
$ws.Range("E45").Value = '  +3.76%  '

# Row 46
$ws.Range("D46").Value = '21.48'
$ws.Range("E46").Value = '  -3.21%  '

# Row 47
$ws.Range("E47").Value = '  +0.11%  '

# Row 48
$ws.Range("D48").Value = '0.0248'
$ws.Range("E48").Value = '  +0.28%  '

# Row 49
$ws.Range("D49").Value = '1.98'
$ws.Range("E49").Value = '  +3.29%  '

# Row 50
$ws.Range("E50").Value = '  +0.50%  '

# Row 51
$ws.Range("D51").Value = '19.27'
$ws.Range("E51").Value = '  +2.36%  '
